$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update LR-pair metrics for F2 -> Gp9 (OldD7) with new TPM-based values
$ws.Cells.Item(2,7).Value = 0.754521
$ws.Cells.Item(2,8).Value = 2.263563
$ws.Cells.Item(2,9).Value = 0.2768403531129761
$ws.Cells.Item(2,10).Value = 0.2768403531129761
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.2089833333333333
$ws.Cells.Item(2,14).Value = 0.62695
$ws.Cells.Item(2,15).Value = 0.04265890843137216
$ws.Cells.Item(2,16).Value = 0.04265890843137217
$ws.Cells.Item(2,17).Value = 0.15768231365
$ws.Cells.Item(2,18).Value = 1.41914082285
$ws.Cells.Item(2,19).Value = 0.01180970727355518
$ws.Cells.Item(2,20).Value = 0.01180970727355518
$ws.Cells.Item(3,7).Value = 0.754521
$ws.Cells.Item(3,8).Value = 2.263563
$ws.Cells.Item(3,9).Value = 0.2768403531129761
$ws.Cells.Item(3,10).Value = 0.2768403531129761
$ws.Cells.Item(3,15).Value = 0.3057562207534381
$ws.Cells.Item(3,16).Value = 0.3057562207534382
$ws.Cells.Item(3,17).Value = 1.130182418494
$ws.Cells.Item(3,18).Value = 10.171641766446
$ws.Cells.Item(3,19).Value = 0.08464566011987089
$ws.Cells.Item(3,20).Value = 0.0846456601198709
$ws.Cells.Item(4,7).Value = 0.754521
$ws.Cells.Item(4,8).Value = 2.263563
$ws.Cells.Item(4,9).Value = 0.2768403531129761
$ws.Cells.Item(4,10).Value = 0.2768403531129761
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.404158
$ws.Cells.Item(4,14).Value = 1.212474
$ws.Cells.Item(4,15).Value = 0.08249911052144433
$ws.Cells.Item(4,16).Value = 0.08249911052144433
$ws.Cells.Item(4,17).Value = 0.3049456983179999
$ws.Cells.Item(4,18).Value = 2.744511284861999
$ws.Cells.Item(4,19).Value = 0.02283908288826309
$ws.Cells.Item(4,20).Value = 0.02283908288826309
$ws.Cells.Item(5,7).Value = 0.754521
$ws.Cells.Item(5,8).Value = 2.263563
$ws.Cells.Item(5,9).Value = 0.2768403531129761
$ws.Cells.Item(5,10).Value = 0.2768403531129761
$ws.Cells.Item(5,13).Value = 2.787915666666667
$ws.Cells.Item(5,14).Value = 8.363747
$ws.Cells.Item(5,15).Value = 0.5690857602937452
$ws.Cells.Item(5,16).Value = 0.5690857602937454
$ws.Cells.Item(5,17).Value = 2.103540916729
$ws.Cells.Item(5,18).Value = 18.931868250561
$ws.Cells.Item(5,19).Value = 0.1575459028312869
$ws.Cells.Item(5,20).Value = 0.1575459028312869
$ws.Cells.Item(6,7).Value = 0.9731926666666667
$ws.Cells.Item(6,9).Value = 0.3570728998754956
$ws.Cells.Item(6,10).Value = 0.3570728998754956
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.2089833333333333
$ws.Cells.Item(6,14).Value = 0.62695
$ws.Cells.Item(6,15).Value = 0.04265890843137216
$ws.Cells.Item(6,16).Value = 0.04265890843137217
$ws.Cells.Item(6,17).Value = 0.2033810474555555
$ws.Cells.Item(6,18).Value = 1.8304294271
$ws.Cells.Item(6,19).Value = 0.01523234013911329
$ws.Cells.Item(6,20).Value = 0.01523234013911329
$ws.Cells.Item(7,7).Value = 0.9731926666666667
$ws.Cells.Item(7,9).Value = 0.3570728998754956
$ws.Cells.Item(7,10).Value = 0.3570728998754956
$ws.Cells.Item(7,15).Value = 0.3057562207534381
$ws.Cells.Item(7,16).Value = 0.3057562207534382
$ws.Cells.Item(7,19).Value = 0.1091772603994023
$ws.Cells.Item(7,20).Value = 0.1091772603994024
$ws.Cells.Item(8,7).Value = 0.9731926666666667
$ws.Cells.Item(8,9).Value = 0.3570728998754956
$ws.Cells.Item(8,10).Value = 0.3570728998754956
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.404158
$ws.Cells.Item(8,14).Value = 1.212474
$ws.Cells.Item(8,15).Value = 0.08249911052144433
$ws.Cells.Item(8,16).Value = 0.08249911052144433
$ws.Cells.Item(8,17).Value = 0.3933236017746666
$ws.Cells.Item(8,18).Value = 3.539912415972
$ws.Cells.Item(8,19).Value = 0.02945819663104114
$ws.Cells.Item(8,20).Value = 0.02945819663104114
$ws.Cells.Item(9,7).Value = 0.9731926666666667
$ws.Cells.Item(9,9).Value = 0.3570728998754956
$ws.Cells.Item(9,10).Value = 0.3570728998754956
$ws.Cells.Item(9,13).Value = 2.787915666666667
$ws.Cells.Item(9,14).Value = 8.363747
$ws.Cells.Item(9,15).Value = 0.5690857602937452
$ws.Cells.Item(9,16).Value = 0.5690857602937454
$ws.Cells.Item(9,17).Value = 2.713179082085111
$ws.Cells.Item(9,18).Value = 24.418611738766
$ws.Cells.Item(9,19).Value = 0.2032051027059388
$ws.Cells.Item(9,20).Value = 0.2032051027059388
$ws.Cells.Item(10,7).Value = 0.7824410000000001
$ws.Cells.Item(10,8).Value = 2.347323
$ws.Cells.Item(10,9).Value = 0.2870844452706686
$ws.Cells.Item(10,10).Value = 0.2870844452706686
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.2089833333333333
$ws.Cells.Item(10,14).Value = 0.62695
$ws.Cells.Item(10,15).Value = 0.04265890843137216
$ws.Cells.Item(10,16).Value = 0.04265890843137217
$ws.Cells.Item(10,17).Value = 0.1635171283166667
$ws.Cells.Item(10,18).Value = 1.47165415485
$ws.Cells.Item(10,19).Value = 0.01224670906287272
$ws.Cells.Item(10,20).Value = 0.01224670906287273
$ws.Cells.Item(11,7).Value = 0.7824410000000001
$ws.Cells.Item(11,8).Value = 2.347323
$ws.Cells.Item(11,9).Value = 0.2870844452706686
$ws.Cells.Item(11,10).Value = 0.2870844452706686
$ws.Cells.Item(11,15).Value = 0.3057562207534381
$ws.Cells.Item(11,16).Value = 0.3057562207534382
$ws.Cells.Item(11,17).Value = 1.172003246707333
$ws.Cells.Item(11,18).Value = 10.548029220366
$ws.Cells.Item(11,19).Value = 0.08777785502305686
$ws.Cells.Item(11,20).Value = 0.08777785502305689
$ws.Cells.Item(12,7).Value = 0.7824410000000001
$ws.Cells.Item(12,8).Value = 2.347323
$ws.Cells.Item(12,9).Value = 0.2870844452706686
$ws.Cells.Item(12,10).Value = 0.2870844452706686
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.404158
$ws.Cells.Item(12,14).Value = 1.212474
$ws.Cells.Item(12,15).Value = 0.08249911052144433
$ws.Cells.Item(12,16).Value = 0.08249911052144433
$ws.Cells.Item(12,17).Value = 0.316229789678
$ws.Cells.Item(12,18).Value = 2.846068107102
$ws.Cells.Item(12,19).Value = 0.02368421137937242
$ws.Cells.Item(12,20).Value = 0.02368421137937243
$ws.Cells.Item(13,7).Value = 0.7824410000000001
$ws.Cells.Item(13,8).Value = 2.347323
$ws.Cells.Item(13,9).Value = 0.2870844452706686
$ws.Cells.Item(13,10).Value = 0.2870844452706686
$ws.Cells.Item(13,13).Value = 2.787915666666667
$ws.Cells.Item(13,14).Value = 8.363747
$ws.Cells.Item(13,15).Value = 0.5690857602937452
$ws.Cells.Item(13,16).Value = 0.5690857602937454
$ws.Cells.Item(13,17).Value = 2.181379522142334
$ws.Cells.Item(13,18).Value = 19.632415699281
$ws.Cells.Item(13,19).Value = 0.1633756698053665
$ws.Cells.Item(13,20).Value = 0.1633756698053666
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.2153186666666667
$ws.Cells.Item(14,8).Value = 0.645956
$ws.Cells.Item(14,9).Value = 0.07900230174085969
$ws.Cells.Item(14,10).Value = 0.07900230174085969
$ws.Cells.Item(14,11).Value = 1
$ws.Cells.Item(14,12).Value = 0.3333333333333333
$ws.Cells.Item(14,13).Value = 0.2089833333333333
$ws.Cells.Item(14,14).Value = 0.62695
$ws.Cells.Item(14,15).Value = 0.04265890843137216
$ws.Cells.Item(14,16).Value = 0.04265890843137217
$ws.Cells.Item(14,17).Value = 0.04499801268888889
$ws.Cells.Item(14,18).Value = 0.4049821142
$ws.Cells.Item(14,19).Value = 0.003370151955830967
$ws.Cells.Item(14,20).Value = 0.003370151955830968
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.2153186666666667
$ws.Cells.Item(15,8).Value = 0.645956
$ws.Cells.Item(15,9).Value = 0.07900230174085969
$ws.Cells.Item(15,10).Value = 0.07900230174085969
$ws.Cells.Item(15,15).Value = 0.3057562207534381
$ws.Cells.Item(15,16).Value = 0.3057562207534382
$ws.Cells.Item(15,17).Value = 0.3225216679724444
$ws.Cells.Item(15,18).Value = 2.902695011751999
$ws.Cells.Item(15,19).Value = 0.02415544521110802
$ws.Cells.Item(15,20).Value = 0.02415544521110803
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.2153186666666667
$ws.Cells.Item(16,8).Value = 0.645956
$ws.Cells.Item(16,9).Value = 0.07900230174085969
$ws.Cells.Item(16,10).Value = 0.07900230174085969
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.404158
$ws.Cells.Item(16,14).Value = 1.212474
$ws.Cells.Item(16,15).Value = 0.08249911052144433
$ws.Cells.Item(16,16).Value = 0.08249911052144433
$ws.Cells.Item(16,17).Value = 0.08702276168266665
$ws.Cells.Item(16,18).Value = 0.7832048551439998
$ws.Cells.Item(16,19).Value = 0.006517619622767677
$ws.Cells.Item(16,20).Value = 0.006517619622767677
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 0.2153186666666667
$ws.Cells.Item(17,8).Value = 0.645956
$ws.Cells.Item(17,9).Value = 0.07900230174085969
$ws.Cells.Item(17,10).Value = 0.07900230174085969
$ws.Cells.Item(17,13).Value = 2.787915666666667
$ws.Cells.Item(17,14).Value = 8.363747
$ws.Cells.Item(17,15).Value = 0.5690857602937452
$ws.Cells.Item(17,16).Value = 0.5690857602937454
$ws.Cells.Item(17,17).Value = 0.6002902841257778
$ws.Cells.Item(17,18).Value = 5.402612557132
$ws.Cells.Item(17,19).Value = 0.04495908495115301
$ws.Cells.Item(17,20).Value = 0.04495908495115302
